$wb = $excel.ActiveWorkbook

# C-column values (Reaction_number) for the 15 new rows (rows 3..17) on each sheet.
# A-column values run 1..15, B-column values run 6..20 (continuing the existing
# row 2 pattern of A=0, B=5).
$sheet1C = @(832, 867, 866, 864, 862, 860, 861, 860, 0, 851, 848, 843, 838, 824, 823)
$sheet2C = @(1017, 935, 936, 935, 939, 936, 929, 925, 0, 913, 912, 922, 918, 912, 914)

$sheetsData = @{
    1 = $sheet1C
    2 = $sheet2C
}

foreach ($sheetIndex in 1, 2) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $cValues = $sheetsData[$sheetIndex]

    for ($i = 0; $i -lt 15; $i++) {
        $row = $i + 3
        $ws.Cells.Item($row, 1).Value = $i + 1
        $ws.Cells.Item($row, 2).Value = $i + 6
        $ws.Cells.Item($row, 3).Value = $cValues[$i]
    }

    # Column A carries the same style as A2 (bordered, bold, centered) - copy it
    # down onto the new rows so the new cells match the existing formatting.
    $ws.Range("A2").Copy()
    $ws.Range("A3:A17").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
